$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive new values for rows 2..25 (bus index 0..23):
# B, C, D, E, I, J, K, L, N
$cols = @("B", "C", "D", "E", "I", "J", "K", "L", "N")

# New vm_pu values per row (row 2 = bus 0 ... row 25 = bus 23), in $cols order.
$data = @(
    @(1.02, 1.013661030557652, 1.019835330644441, 1.015557008224913, 1.025415332750196, 1.018896720792565, 1.022677400990017, 1.018411811045556, 1.010454810447489),
    @(1.02, 1.014442800291442, 1.020383181488358, 1.016214634519904, 1.025497789374017, 1.019313614728725, 1.023032071704415, 1.01887503543619, 1.010591710862483),
    @(1.02, 1.014948865474225, 1.020737272801257, 1.016640763650113, 1.025549131156304, 1.019582934879522, 1.023260458189376, 1.019174692318271, 1.010680146969364),
    @(1.02, 1.015161663175033, 1.020886033737131, 1.016820050202678, 1.025570232203532, 1.019696050731673, 1.023356204686031, 1.019300647317738, 1.010717289681996),
    @(1.02, 1.015197395599602, 1.020911005481954, 1.016850161464508, 1.025573746805007, 1.0197150371027, 1.02337226522854, 1.019321794454557, 1.010723523991773),
    @(1.02, 1.014951708701672, 1.020739260944086, 1.016643158730899, 1.025549415008914, 1.01958444675876, 1.023261738610677, 1.019176375418632, 1.010680643413501),
    @(1.02, 1.013925189131205, 1.020020561940624, 1.015779130480726, 1.025443615543995, 1.019037701888826, 1.022797492298987, 1.01856837559704, 1.010501106856424),
    @(1.02, 1.012118005098458, 1.018751116018757, 1.014261292929665, 1.025241812318482, 1.018070976272769, 1.021971008005278, 1.017496455018905, 1.010183631363023),
    @(1.02, 1.010914452184691, 1.01790292366921, 1.013252673597427, 1.025097005657123, 1.017424365637327, 1.021414452003011, 1.01678156426223, 1.009971264076356),
    @(1.02, 1.010393616891236, 1.017535224076091, 1.012816731464153, 1.02503187981917, 1.017143889659535, 1.021172157066827, 1.016471962065505, 1.009879142679002),
    @(1.02, 1.010200204054485, 1.017398582101054, 1.012654924865399, 1.025007326099229, 1.0170396363228, 1.02108196386386, 1.016356956081656, 1.009844900389642),
    @(1.02, 1.010241689536493, 1.017427895038889, 1.012689627389588, 1.025012609369301, 1.017062002255684, 1.021101319368608, 1.016381625499934, 1.009852246568879),
    @(1.02, 1.010377628324099, 1.017523930459775, 1.012803353985495, 1.025029857601871, 1.017135273508099, 1.021164705618328, 1.016462455747382, 1.009876312693625),
    @(1.02, 1.010461391190391, 1.017583092903545, 1.012873440884179, 1.025040436725864, 1.017180408837953, 1.021203734309956, 1.016512257188889, 1.009891137421685),
    @(1.02, 1.010949025019345, 1.017927317872402, 1.013281622599607, 1.025101276870539, 1.017442969718419, 1.021430505015355, 1.016802110622834, 1.00997737443584),
    @(1.02, 1.011254989423166, 1.018143128369629, 1.013537878936036, 1.025138792152572, 1.017607537241465, 1.021572404709313, 1.016983915783017, 1.010031424813423),
    @(1.02, 1.01143348317259, 1.018268965667891, 1.013687425655587, 1.025160440270906, 1.017703479386191, 1.021655046618734, 1.017089954751256, 1.010062935535642),
    @(1.02, 1.011494349941871, 1.018311865862062, 1.01373843017598, 1.025167782012879, 1.017736185110064, 1.021683203990278, 1.017126110420179, 1.010073677157729),
    @(1.02, 1.01122215923331, 1.018119978203506, 1.013510377104725, 1.025134791304866, 1.017589885586466, 1.021557193213555, 1.016964410306238, 1.010025627365578),
    @(1.02, 1.010337596383432, 1.017495652115669, 1.012769860961275, 1.02502478844191, 1.017113698928451, 1.021146045292182, 1.016438653405905, 1.009869226486672),
    @(1.02, 1.009781718599878, 1.017102756007106, 1.012304974300607, 1.024953525090093, 1.016813885208278, 1.020886418205273, 1.016108055546364, 1.009770750687525),
    @(1.02, 1.010076372585484, 1.017311070847631, 1.012551351993743, 1.02499150187137, 1.016972861111811, 1.021024157259977, 1.016283314403497, 1.009822967715737),
    @(1.02, 1.011236993690804, 1.018130438894999, 1.013522803761651, 1.025136599838585, 1.017597861755893, 1.021564067022545, 1.016973224006935, 1.010028247032046),
    @(1.02, 1.012584995629084, 1.019079641433518, 1.014653122195489, 1.025295798268616, 1.018321279163167, 1.022185662170961, 1.017773627636935, 1.010265834852989)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $vals[$j]
    }
}
